# Hortaliza, Vega Monumental Concepción - Brócoli
# Weekly data refresh: insert the latest week's two rows (Primera / Segunda
# quality) at the top of the data block (rows 137-138), pushing the
# pre-existing rows down by two (the oldest pair now lands at 153-154).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 137, shifting all the
# existing data (old rows 137-152) down to 139-154.
$ws.Rows("137:138").Insert()

# New "Primera" quality record for the latest week.
$row137 = @(11, "Vega Monumental Concepción", "Bíobío", 44474, 8, 100112023, `
    "Brócoli", "Sin especificar", "Primera", 1000, 600, 700, 650, `
    "$/unidad", "Región Metropolitana", 650, 1, "Hortaliza")

# New "Segunda" quality record for the latest week.
$row138 = @(11, "Vega Monumental Concepción", "Bíobío", 44474, 8, 100112023, `
    "Brócoli", "Sin especificar", "Segunda", 500, 500, 500, 500, `
    "$/unidad", "Región Metropolitana", 500, 1, "Hortaliza")

for ($c = 1; $c -le 18; $c++) {
    $ws.Cells.Item(137, $c).Value = $row137[$c - 1]
    $ws.Cells.Item(138, $c).Value = $row138[$c - 1]
}
